{"js": "// The document is a single-column table where each row holds one stat\n// value. Apply the value replacements from the diff by row index, and\n// collapse the three multi-run \"detail\" rows (11/88/3 + tab-separated\n// figures) down to a single run holding the corresponding summary value.\n\nconst table = context.document.body.tables.getFirst();\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst updates = [\n  [0, \"0M\"],\n  [1, \"0M\"],\n  [2, \"0M\"],\n  [3, \"714\"],\n  [5, \"0.09867\"],\n  [6, \"0.01876\"],\n  [7, \"0.00576\"],\n  [8, \"0.04055\"],\n  [9, \"0.05590\"],\n  [10, \"0.06801\"],\n  [11, \"5.05377\"],\n  [43, \"95.7\"],\n  [44, \"5.05\"],\n  [45, \"117\"],\n];\n\nfor (const [rowIndex, newText] of updates) {\n  table.getCell(rowIndex, 0).value = newText;\n}\n\nawait context.sync();\n", "ps1": "# The document is a single-column table where each row holds one stat\n# value. Apply the value replacements from the diff by (1-based) row\n# index, and collapse the three multi-run \"detail\" rows (11/88/3 +\n# tab-separated figures) down to a single run holding the corresponding\n# summary value.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text  = \"0M\"\n$t.Cell(2,1).Range.Text  = \"0M\"\n$t.Cell(3,1).Range.Text  = \"0M\"\n$t.Cell(4,1).Range.Text  = \"714\"\n$t.Cell(6,1).Range.Text  = \"0.09867\"\n$t.Cell(7,1).Range.Text  = \"0.01876\"\n$t.Cell(8,1).Range.Text  = \"0.00576\"\n$t.Cell(9,1).Range.Text  = \"0.04055\"\n$t.Cell(10,1).Range.Text = \"0.05590\"\n$t.Cell(11,1).Range.Text = \"0.06801\"\n$t.Cell(12,1).Range.Text = \"5.05377\"\n$t.Cell(44,1).Range.Text = \"95.7\"\n$t.Cell(45,1).Range.Text = \"5.05\"\n$t.Cell(46,1).Range.Text = \"117\"\n"}
